# paises.xlsx - refresh COVID-19 country stats snapshot
# (commit: "Update countries & provincias Spain")
#
# Sheet "Pais" is sorted descending by "Casos totales" (col B). A newer data
# pull (footer timestamp "...14:58" -> "...16:15") bumped several countries'
# case counts. That re-shuffles a handful of table rows relative to their
# neighbours (countries overtaking others in the ranking) and refreshes other
# countries' figures in place. Each affected row below is written with its
# final country name (col A) and final Casos totales / Nuevos casos / Casos
# activos / Recuperados / Casos criticos / Muertes hoy / Muertes (cols B-H).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer timestamp: "...a las 14:58" -> "...a las 16:15"
$ws.Range("A1").Value = "Datos actualizados a 17 de Junio de 2020 a las 16:15"

# Estados Unidos (4), Brasil (5): updated case counts
$ws.Cells.Item(4, 2).Value = 2209260
$ws.Cells.Item(4, 3).Value = 860
$ws.Cells.Item(4, 5).Value = 1186963
$ws.Cells.Item(4, 7).Value = 29
$ws.Cells.Item(4, 8).Value = 119161
$ws.Cells.Item(5, 2).Value = 929149
$ws.Cells.Item(5, 3).Value = 315
$ws.Cells.Item(5, 5).Value = 406318
$ws.Cells.Item(5, 7).Value = 11
$ws.Cells.Item(5, 8).Value = 45467

# India (7): updated case counts
$ws.Cells.Item(7, 2).Value = 356372
$ws.Cells.Item(7, 3).Value = 2211
$ws.Cells.Item(7, 5).Value = 156700
$ws.Cells.Item(7, 7).Value = 33
$ws.Cells.Item(7, 8).Value = 11954

# Alemania (13): updated case counts
$ws.Cells.Item(13, 2).Value = 188544
$ws.Cells.Item(13, 3).Value = 162
$ws.Cells.Item(13, 5).Value = 6034

# Suecia (28): updated case counts
$ws.Cells.Item(28, 2).Value = 54562
$ws.Cells.Item(28, 3).Value = 1239
$ws.Cells.Item(28, 7).Value = 102
$ws.Cells.Item(28, 8).Value = 5041

# Re-rank: Ghana overtakes Serbia/Dinamarca/Moldavia/Corea del Sur,
# shifting rows 57-61 down one rank (new countries + refreshed figures)
$ws.Cells.Item(57, 1).Value = "Ghana"
$ws.Cells.Item(57, 2).Value = 12590
$ws.Cells.Item(57, 3).Value = 397
$ws.Cells.Item(57, 4).Value = 4410
$ws.Cells.Item(57, 5).Value = 8114
$ws.Cells.Item(57, 7).Value = 8
$ws.Cells.Item(57, 8).Value = 66
$ws.Cells.Item(58, 1).Value = "Serbia"
$ws.Cells.Item(58, 2).Value = 12522
$ws.Cells.Item(58, 3).Value = 96
$ws.Cells.Item(58, 4).Value = 11511
$ws.Cells.Item(58, 5).Value = 754
$ws.Cells.Item(58, 7).Value = 1
$ws.Cells.Item(58, 8).Value = 257
$ws.Cells.Item(59, 1).Value = "Dinamarca"
$ws.Cells.Item(59, 2).Value = 12294
$ws.Cells.Item(59, 3).Value = 44
$ws.Cells.Item(59, 4).Value = 11185
$ws.Cells.Item(59, 5).Value = 511
$ws.Cells.Item(59, 7).Value = 0
$ws.Cells.Item(59, 8).Value = 598
$ws.Cells.Item(60, 1).Value = "Moldavia"
$ws.Cells.Item(60, 2).Value = 12254
$ws.Cells.Item(60, 3).Value = 0
$ws.Cells.Item(60, 4).Value = 7077
$ws.Cells.Item(60, 5).Value = 4750
$ws.Cells.Item(60, 7).Value = 4
$ws.Cells.Item(60, 8).Value = 427
$ws.Cells.Item(61, 1).Value = "Corea del Sur"
$ws.Cells.Item(61, 2).Value = 12198
$ws.Cells.Item(61, 3).Value = 43
$ws.Cells.Item(61, 4).Value = 10774
$ws.Cells.Item(61, 5).Value = 1145
$ws.Cells.Item(61, 7).Value = 1
$ws.Cells.Item(61, 8).Value = 279

# Noruega (69): updated case counts
$ws.Cells.Item(69, 2).Value = 8680
$ws.Cells.Item(69, 3).Value = 20
$ws.Cells.Item(69, 5).Value = 299

# Cuba (99): updated case counts
$ws.Cells.Item(99, 2).Value = 2280
$ws.Cells.Item(99, 3).Value = 7
$ws.Cells.Item(99, 4).Value = 1999
$ws.Cells.Item(99, 5).Value = 197

# Sri Lanka (104), Mali (105): updated case counts
$ws.Cells.Item(104, 2).Value = 1924
$ws.Cells.Item(104, 3).Value = 9
$ws.Cells.Item(104, 5).Value = 516
$ws.Cells.Item(105, 2).Value = 1890
$ws.Cells.Item(105, 3).Value = 5
$ws.Cells.Item(105, 4).Value = 1168
$ws.Cells.Item(105, 5).Value = 615
$ws.Cells.Item(105, 7).Value = 3
$ws.Cells.Item(105, 8).Value = 107

# Islandia (107): updated case counts
$ws.Cells.Item(107, 2).Value = 1815
$ws.Cells.Item(107, 3).Value = 3
$ws.Cells.Item(107, 4).Value = 1797
$ws.Cells.Item(107, 5).Value = 8

# Re-rank: Guayana Francesa overtakes Nueva Zelanda/Eslovenia/Guinea-Bisau/Libano,
# shifting rows 114-118 down one rank (new country + refreshed figures)
$ws.Cells.Item(114, 1).Value = "Guayana Francesa"
$ws.Cells.Item(114, 2).Value = 1554
$ws.Cells.Item(114, 3).Value = 133
$ws.Cells.Item(114, 4).Value = 663
$ws.Cells.Item(114, 5).Value = 886
$ws.Cells.Item(114, 8).Value = 5
$ws.Cells.Item(115, 1).Value = "Nueva Zelanda"
$ws.Cells.Item(115, 2).Value = 1506
$ws.Cells.Item(115, 3).Value = 0
$ws.Cells.Item(115, 4).Value = 1482
$ws.Cells.Item(115, 5).Value = 2
$ws.Cells.Item(115, 8).Value = 22
$ws.Cells.Item(116, 1).Value = "Eslovenia"
$ws.Cells.Item(116, 2).Value = 1503
$ws.Cells.Item(116, 3).Value = 4
$ws.Cells.Item(116, 4).Value = 1359
$ws.Cells.Item(116, 5).Value = 35
$ws.Cells.Item(116, 8).Value = 109
$ws.Cells.Item(117, 1).Value = "Guinea-Bisau"
$ws.Cells.Item(117, 2).Value = 1492
$ws.Cells.Item(117, 4).Value = 153
$ws.Cells.Item(117, 5).Value = 1324
$ws.Cells.Item(117, 8).Value = 15
$ws.Cells.Item(118, 1).Value = "Libano"
$ws.Cells.Item(118, 2).Value = 1489
$ws.Cells.Item(118, 3).Value = 16
$ws.Cells.Item(118, 4).Value = 907
$ws.Cells.Item(118, 5).Value = 550
$ws.Cells.Item(118, 8).Value = 32

# Tunez (123): updated case counts
$ws.Cells.Item(123, 2).Value = 1128
$ws.Cells.Item(123, 3).Value = 3
$ws.Cells.Item(123, 4).Value = 1004
$ws.Cells.Item(123, 7).Value = 1
$ws.Cells.Item(123, 8).Value = 50

# Cabo Verde (136): updated case counts
$ws.Cells.Item(136, 2).Value = 791
$ws.Cells.Item(136, 3).Value = 10
$ws.Cells.Item(136, 5).Value = 430

# Birmania (160): updated case counts
$ws.Cells.Item(160, 4).Value = 185
$ws.Cells.Item(160, 5).Value = 71

# Swap: Islas Malvinas <-> Groenlandia (identical totals, order only)
$ws.Cells.Item(206, 1).Value = "Islas Malvinas"
$ws.Cells.Item(207, 1).Value = "Groenlandia"

# Swap: Montserrat <-> Seychelles (ranks tied on Casos totales)
$ws.Cells.Item(210, 1).Value = "Montserrat"
$ws.Cells.Item(210, 4).Value = 10
$ws.Cells.Item(210, 8).Value = 1
$ws.Cells.Item(211, 1).Value = "Seychelles"
$ws.Cells.Item(211, 4).Value = 11
$ws.Cells.Item(211, 8).Value = 0

# Swap: Islas Virgenes Britanicas <-> Papua Nueva Guinea (ranks tied)
$ws.Cells.Item(213, 1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(213, 4).Value = 7
$ws.Cells.Item(213, 8).Value = 1
$ws.Cells.Item(214, 1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(214, 4).Value = 8
$ws.Cells.Item(214, 8).Value = 0
